$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear C2 and E2 (forecast values for the first row are removed)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: clear C3, update E3
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = 2.807231216534301

# Row 4: update C4, E4
$ws.Range("C4").Value = -0.9140166223623458
$ws.Range("E4").Value = 1.821983295885099

# Row 8: update C8
$ws.Range("C8").Value = -1.479696720105184

# Row 9: update E9
$ws.Range("E9").Value = -0.6155071485167807

# Row 11: update C11, E11
$ws.Range("C11").Value = 2.192778679161966
$ws.Range("E11").Value = -0.5835597102573087

# Row 12: update C12
$ws.Range("C12").Value = 3.408364488606752

# Row 13: update E13
$ws.Range("E13").Value = 3.056075254340018

# Row 15: update C15, E15
$ws.Range("C15").Value = 1.666553973046025
$ws.Range("E15").Value = -1.376301649685407

# Row 16: update C16, E16
$ws.Range("C16").Value = 1.879266440112781
$ws.Range("E16").Value = -0.5015683214423916

# Row 17: update C17, E17
$ws.Range("C17").Value = -2.620683231370935
$ws.Range("E17").Value = -3.531225750971467

# Row 18: update C18
$ws.Range("C18").Value = -3.036556262700263
